$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I18").Value = 0.2019572604821628
$ws.Range("J18").Value = 0.1159468389252036
$ws.Range("K18").Value = 0.07503651310250715
$ws.Range("L18").Value = 2.189568379381487
